$d = $word.ActiveDocument
$xml = $d.Content.WordOpenXML
Write-Output $xml
